# Add "2022-Q4" sheet (new quarter data) and update the "总计" (totals) sheet
# with the new quarter's row, per commit "feat: add 2022-Q4 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet right after "总计" (i.e. before "2022-Q3"),
#    matching the target sheet order: 总计, 2022-Q4, 2022-Q3, 2022-Q2, ...
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item(1)
$added = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetTotal)
$added.Name = "2022-Q4"

# Re-resolve worksheet references AFTER the Add() call - references captured
# before a sheet insertion can end up stale for cross-sheet operations.
$sheetTotal = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Item(2)
$q3src = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 2. Build the "2022-Q4" sheet. Column layout mirrors the other quarter
#    sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
#    Pull header + index-column cell formatting from the existing "2022-Q3"
#    sheet (single-cell Copy/PasteSpecial -4122 = xlPasteFormats) so the
#    bordered/bold style used there (s="2") is reproduced exactly.
# ---------------------------------------------------------------------------
$headers = @{
    "B" = "基金代码";
    "C" = "基金名称";
    "D" = "基金规模";
    "E" = "股票总仓位";
    "F" = "仓位占比";
    "G" = "持有市值(亿元)";
    "H" = "仓位排名";
}

foreach ($col in @("B","C","D","E","F","G","H")) {
    $ref = $col + "1"
    $q3src.Range($ref).Copy()
    $q4.Range($ref).PasteSpecial(-4122)
    $q4.Range($ref).Value = $headers[$col]
}

$rows = @(
    @{ idx = 0; code = "167506"; name = "安信深圳科技指数（LOF）A"; scale = "0.87"; pos = "94.03"; ratio = "3.94"; mv = "0.0343"; rank = 7 },
    @{ idx = 1; code = "167507"; name = "安信深圳科技指数（LOF）C"; scale = "0.33"; pos = "94.03"; ratio = "3.94"; mv = "0.0130"; rank = 7 },
    @{ idx = 2; code = "006611"; name = "人保中证500指数";         scale = "0.43"; pos = "93.30"; ratio = "0.47"; mv = "0.0020"; rank = 5 },
    @{ idx = 3; code = "510570"; name = "兴业中证500ETF";          scale = "0.09"; pos = "97.48"; ratio = "0.58"; mv = "0.0005"; rank = 3 }
)

$r = 2
foreach ($row in $rows) {
    # A column: numeric index, with the same bordered/bold style as the header row.
    $q3src.Range("A2").Copy()
    $q4.Range("A" + $r).PasteSpecial(-4122)
    $q4.Range("A" + $r).Value = $row.idx

    # B column (fund code): force text so codes like "006611" keep their
    # leading zero instead of being parsed as a number.
    $q4.Range("B" + $r).NumberFormat = "@"
    $q4.Range("B" + $r).Value = $row.code
    $q4.Range("C" + $r).Copy()
    $q4.Range("B" + $r).PasteSpecial(-4122)

    # C column (fund name): plain Chinese text, never parsed as a number.
    $q4.Range("C" + $r).Value = $row.name

    # D/E/F/G columns: decimal-looking values stored as literal text in the
    # source data, so force text the same way, then strip the Text number
    # format back off (paste formats from the already-plain C column) so the
    # cell ends up with no explicit style, matching the sibling sheets.
    $q4.Range("D" + $r).NumberFormat = "@"
    $q4.Range("D" + $r).Value = $row.scale
    $q4.Range("C" + $r).Copy()
    $q4.Range("D" + $r).PasteSpecial(-4122)

    $q4.Range("E" + $r).NumberFormat = "@"
    $q4.Range("E" + $r).Value = $row.pos
    $q4.Range("C" + $r).Copy()
    $q4.Range("E" + $r).PasteSpecial(-4122)

    $q4.Range("F" + $r).NumberFormat = "@"
    $q4.Range("F" + $r).Value = $row.ratio
    $q4.Range("C" + $r).Copy()
    $q4.Range("F" + $r).PasteSpecial(-4122)

    $q4.Range("G" + $r).NumberFormat = "@"
    $q4.Range("G" + $r).Value = $row.mv
    $q4.Range("C" + $r).Copy()
    $q4.Range("G" + $r).PasteSpecial(-4122)

    # H column (rank): real number.
    $q4.Range("H" + $r).Value = $row.rank

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Update "总计" sheet: insert a new row right under the header for the
#    2022-Q4 totals, pushing the existing quarters down (the row index
#    column keeps counting 0,1,2,... so no other value needs to change).
# ---------------------------------------------------------------------------
$sheetTotal.Rows("2:2").Insert()

# Re-apply the bordered/bold index-column style to the new A2 (row insert
# does not carry over the bottom-row style for that column), and scrub the
# stray formatting the insert leaves on B2:D2 so they stay at the default
# style like every other data row in the table.
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)
$sheetTotal.Range("B2:D2").ClearFormats()

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 4
$sheetTotal.Range("D2").Value = 0.05
